# The paragraph ending in "...Tuer, ou être tuée." originally reads (run by
# run):
#   "Elle" | [_GoBack bookmark] | " avait fait ce qu'il fallait. Rien
#   d'autre que se défendre. " | "Tuer, ou être tuée."
#
# The edit moves the stray "_GoBack" bookmark (Word's last-edit-position
# marker) from right after "Elle" to right after "Tuer," — i.e. between
# "Tuer," and " ou être tuée." — without altering the visible text at all.
#
# Net effect once runs are re-split around the bookmark's new position:
#   "Elle" | " avait fait ce qu'il fallait. Rien d'autre que se défendre. "
#   | "Tuer," | [_GoBack bookmark] | " ou être tuée."

$d = $word.ActiveDocument

# 1. Remove the bookmark from its current (stale) location, if present.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Cut the non-italic " avait fait ce qu'il fallait. Rien d'autre que se
#    défendre. " run, currently sitting right after the old bookmark spot.
$middle = $d.Content
$middle.Find.Execute(" avait fait ce qu’il fallait. Rien d’autre que se défendre. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$middle.Cut()

# 3. Re-insert it immediately after "Elle" (where the bookmark used to be),
#    restoring the original reading order.
$elleAnchor = $d.Content
$elleAnchor.Find.Execute("pour les conquérir. Elle", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $d.Range($elleAnchor.End, $elleAnchor.End)
$insertPoint.Paste()

# 4. Re-create the "_GoBack" bookmark, collapsed, right after "Tuer," —
#    this is what forces "Tuer," and " ou être tuée." apart into two runs.
$tuerComma = $d.Content
$tuerComma.Find.Execute("Tuer,", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkSpot = $d.Range($tuerComma.End, $tuerComma.End)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null
